$d = $word.ActiveDocument

$pairs = @(
    @("584×5=2920", "835×9=7515"),
    @("299×9=2691", "822×5=4110"),
    @("488×7=3416", "340×6=2040"),
    @("293×5=1465", "216×6=1296"),
    @("208×7=1456", "369×6=2214"),
    @("365×9=3285", "146×8=1168"),
    @("236×9=2124", "749×9=6741"),
    @("529×4=2116", "971×2=1942"),
    @("398×7=2786", "556×6=3336"),
    @("605×7=4235", "474×6=2844"),
    @("740×5=3700", "751×3=2253"),
    @("182×4=728",  "411×3=1233"),
    @("163×6=978",  "386×3=1158"),
    @("898×6=5388", "248×9=2232"),
    @("361×5=1805", "357×6=2142"),
    @("161×9=1449", "601×6=3606"),
    @("666×4=2664", "823×3=2469"),
    @("864×4=3456", "975×9=8775"),
    @("795×6=4770", "764×3=2292"),
    @("558×8=4464", "781×3=2343"),
    @("873×9=7857", "439×2=878"),
    @("270×2=540",  "678×4=2712"),
    @("445×6=2670", "120×8=960"),
    @("432×8=3456", "602×3=1806"),
    @("916×8=7328", "713×6=4278")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
